# Apply updated crypto price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'71.633.47"
$ws.Range("E2").Value = "  +3.24%  "

$ws.Range("D3").Value = "'3.639.82"
$ws.Range("E3").Value = "  +6.94%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").Value = "'589.66"
$ws.Range("E5").Value = "  +0.26%  "

$ws.Range("D6").Value = "'180.30"
$ws.Range("E6").Value = "  -0.65%  "

$ws.Range("D7").Value = "'3.631.56"
$ws.Range("E7").Value = "  +6.90%  "

$ws.Range("D8").Value = "'0.615"
$ws.Range("E8").Value = "  +2.57%  "

$ws.Range("E9").Value = "  +0.07%  "

$ws.Range("E10").Value = "  -0.19%  "

$ws.Range("D11").Value = "'0.607"
$ws.Range("E11").Value = "  +2.16%  "

$ws.Range("D12").Value = "'49.72"
$ws.Range("E12").Value = "  +2.45%  "

$ws.Range("D13").Value = "'0.0000286"

$ws.Range("D14").Value = "'682.29"
$ws.Range("E14").Value = "  -0.52%  "

$ws.Range("D15").Value = "'4.224.86"
$ws.Range("E15").Value = "  +7.03%  "

$ws.Range("D16").Value = "'8.99"
$ws.Range("E16").Value = "  +3.35%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "'71.749.15"
$ws.Range("E17").Value = "  +3.40%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "'3.618.36"
$ws.Range("E18").Value = "  +7.12%  "

$ws.Range("E19").Value = "  +1.85%  "

$ws.Range("D20").Value = "'18.30"
$ws.Range("E20").Value = "  +3.12%  "

$ws.Range("D21").Value = "'11.60"
$ws.Range("E21").Value = "  +2.30%  "

$ws.Range("E22").Value = "  +2.48%  "

$ws.Range("E23").Value = "  +8.97%  "

$ws.Range("D24").Value = "'17.77"
$ws.Range("E24").Value = "  +2.37%  "

$ws.Range("D25").Value = "'103.29"
$ws.Range("E25").Value = "  -0.11%  "

$ws.Range("D26").Value = "'4.02"
$ws.Range("E26").Value = "  +1.69%  "

$ws.Range("E27").Value = "  +3.86%  "

$ws.Range("E28").Value = "  +2.69%  "

$ws.Range("D29").Value = "'35.09"
$ws.Range("E29").Value = "  +3.22%  "

$ws.Range("D30").Value = "'9.19"
$ws.Range("E30").Value = "  +4.03%  "

$ws.Range("D31").Value = "'7.25"
$ws.Range("E31").Value = "  +3.80%  "

$ws.Range("D32").Value = "'4.14"
$ws.Range("E32").Value = "  +13.91%  "

$ws.Range("D33").Value = "'575.12"
$ws.Range("E33").Value = "  +3.05%  "

$ws.Range("D34").Value = "'11.34"
$ws.Range("E34").Value = "  +1.57%  "

$ws.Range("D35").Value = "'0.110"
$ws.Range("E35").Value = "  +2.44%  "

$ws.Range("D36").Value = "'59.44"
$ws.Range("E36").Value = "  +1.45%  "

$ws.Range("E37").Value = "  +0.03%  "

$ws.Range("D38").Value = "'3.671.38"
$ws.Range("E38").Value = "  +0.28%  "

$ws.Range("D39").Value = "'0.143"
$ws.Range("E39").Value = "  +0.54%  "

$ws.Range("D40").Value = "'35.62"
$ws.Range("E40").Value = "  -1.03%  "

$ws.Range("E41").Value = "  +4.07%  "

$ws.Range("D42").Value = "'0.0474"
$ws.Range("E42").Value = "  +10.87%  "

$ws.Range("E43").Value = "  +4.02%  "

$ws.Range("E44").Value = "  +2.47%  "

$ws.Range("E45").Value = "  +2.29%  "

$ws.Range("D46").Value = "'3.39"
$ws.Range("E46").Value = "  +1.42%  "

$ws.Range("D47").Value = "'2.81"
$ws.Range("E47").Value = "  +4.91%  "

$ws.Range("E48").Value = "  +2.48%  "

$ws.Range("E49").Value = "  +3.10%  "

$ws.Range("D50").Value = "'0.998"
$ws.Range("E50").Value = "  -0.24%  "

$ws.Range("B51").Value = "LidoDAOToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D51").Value = "'3.01"
$ws.Range("E51").Value = "  +13.48%  "
